$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Job_Description"
$ws.Range("B1").Value = "Total_Years_Min_Exp"
$ws.Range("C1").Value = "Total_Years_Max_Exp"

$headers = $ws.Range("A1:C1")
$headers.Font.Bold = $true
$headers.Borders.LineStyle = 1
$headers.HorizontalAlignment = -4108
$headers.VerticalAlignment = -4160

# Data row
$ws.Range("A2").Value = "Please find the Job Description (JD) below for your reference.`nJob Role: ML Engineer/Data Scientist`nJob Description:`n• Develop, train, test, and deploy machine learning models across various fields including computer vision, LLMs, and with tabular and time series data.`n• Strong experience in Python, FastAPI, Flask`n• Strong experience in SQL`n• Strong experience Design Pattern/algorithms and data structures`n• Familiarity with OOPS, Design Pattern/algorithms and data structures`n• Familiarity with continuous integration, deployment, and automated build processes for scalable application delivery using Docker/Kubernetes`n• Practical knowledge of one or more major cloud platforms (e.g. Azure, AWS, or GCP).`n• Excellent written and verbal communication skills in English.`n• Experiment with novel deep learning-based technologies such as self-supervised learning and generative AI. `n• Work directly with customer data and set up data pipelines to collect, curate, transform, and version data. `n• Participate in the collection, analysis, interpretation, and output of large amounts of data using advanced AI techniques like deep learning, NLP, and computer vision good foundational experience in PyTorch / Tensorflow.`n• Work within the global corporate Artificial Intelligence division, which addresses real business challenges and opportunities across multiple countries.`n• Collaborate across different business and corporate functions in an international team composed of Project Managers, Data Scientists, Data and Software Engineers within the Artificial Intelligence team and others in the Global AI team`n Qualifications:`n• Bachelor’s degree or master’s degree in data science, Computational Statistics/Mathematics, Computer Science, or related field`n• Fluent English"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 5
